# TouchGFX texts.xlsx - add LCD brightness percentage text entries
# (informacja o jasnosci ekranu na LCD - do dopracowania)

$wb = $excel.ActiveWorkbook

# --- Sheet "Typography": wildcard characters for the numeric font used to
#     render the brightness percentage (row 4, column G = "Wildcard Characters")
$tsTypo = $wb.Worksheets.Item("Typography")
$tsTypo.Range("G4").Value = "0123456789-"

# --- Sheet "Translation": two new single-use text rows driving the new
#     LCD-brightness display text ("<value> %") and its numeric placeholder ("0")
$tsTrans = $wb.Worksheets.Item("Translation")

# Row 6: SingleUseId3 / Typography_00 / Center / LTR / "<value> %"
$tsTrans.Range("B6").Value = "SingleUseId3"
$tsTrans.Range("C6").Value = "Typography_00"
$tsTrans.Range("D6").Value = "Center"
$tsTrans.Range("E6").Value = "LTR"
$tsTrans.Range("F6").Value = "<value> %"

# Row 7: SingleUseId4 / Typography_00 / Center / LTR / "0"
$tsTrans.Range("B7").Value = "SingleUseId4"
$tsTrans.Range("C7").Value = "Typography_00"
$tsTrans.Range("D7").Value = "Center"
$tsTrans.Range("E7").Value = "LTR"

# "0" must stay a plain text label (shared string), not be coerced to a
# number, and must not pick up a quote-prefix / text-number-format style.
$f7 = $tsTrans.Range("F7")
$f7.NumberFormat = "@"
$f7.Value = "0"
$f7.Style = "Normal"
